$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "department" value for the English course row changed from
# "FACULTY OF ENGLISH" to "English".
$ws.Range("C2").Value = "English"
